$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 (DrugUse StartDate block) currently holds the values that belong
# on row 16 (AlcoholUse frequencyAndQuantity.period). Move B13:J13 -> B16:J16
# and blank out B13:J13.

$values = $ws.Range("B13:J13").Value2

$ws.Range("B16:J16").Value2 = $values
$ws.Range("B13:J13").ClearContents()
